$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching the style/formatting of the existing header (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I and J (rows 2-9)
$data = @{
    2 = @(8, 8)
    3 = @(9, 9)
    4 = @(8, 8)
    5 = @(8, 9)
    6 = @(8, 9)
    7 = @(8, 9)
    8 = @(9, 9)
    9 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
